$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$row = 22
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $row++
}

# Match the saved view state: the next empty row (31) is selected as an
# entire-row selection, as would happen after selecting the row below the
# newly pasted/entered data.
$ws.Activate()
[void]$ws.Rows("31:1048576").Select()

# A page setup was defined for the sheet (portrait orientation) when the
# file was resaved.
$ws.PageSetup.Orientation = 1
